$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")
$cell = $ws.Range("B7")
$cell.NumberFormat = "@"
$cell.Value = "123"
